$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '30.162.30'
$ws.Range('E2').Value = '  -3.32%  '
$ws.Range('D3').Value = '1.862.93'
$ws.Range('E3').Value = '  -4.03%  '
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '233.68'
$ws.Range('E5').Value = '  -3.49%  '
$ws.Range('D6').Value = '0.9998'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').Value = '0.4657'
$ws.Range('E7').Value = '  -3.06%  '
$ws.Range('D8').Value = '0.2826'
$ws.Range('E8').Value = '  -3.04%  '
$ws.Range('D9').Value = '0.06542'
$ws.Range('E9').Value = '  -3.63%  '
$ws.Range('D10').Value = '20.06'
$ws.Range('E10').Value = '  -0.76%  '
$ws.Range('D11').Value = '0.07807'
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').Value = '96.14'
$ws.Range('E12').Value = '  -7.85%  '
$ws.Range('D13').Value = '1.857.35'
$ws.Range('D14').Value = '5.122'
$ws.Range('E14').Value = '  -3.47%  '
$ws.Range('D15').Value = '0.6692'
$ws.Range('E15').Value = '  -4.12%  '
$ws.Range('D16').Value = '280.98'
$ws.Range('E16').Value = '  -5.26%  '
$ws.Range('D17').Value = '30.185.43'
$ws.Range('E17').Value = '  -3.23%  '
$ws.Range('D18').Value = '0.9996'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').Value = '5.460'
$ws.Range('E19').Value = '  -1.89%  '
$ws.Range('E20').Value = '  -2.95%  '
$ws.Range('D21').Value = '2.102.37'
$ws.Range('E21').Value = '  -4.76%  '
$ws.Range('D22').Value = '0.000007239'
$ws.Range('E22').Value = '  -4.94%  '
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').Value = '6.146'
$ws.Range('E24').Value = '  -4.55%  '
$ws.Range('D25').Value = '9.319'
$ws.Range('E25').Value = '  -2.66%  '
$ws.Range('D26').Value = '165.23'
$ws.Range('E26').Value = '  -2.30%  '
$ws.Range('D27').Value = '18.89'
$ws.Range('E27').Value = '  -4.86%  '
$ws.Range('D28').Value = '1.905'
$ws.Range('E28').Value = '  -9.56%  '
$ws.Range('D29').Value = '1.342'
$ws.Range('E29').Value = '  -4.08%  '
$ws.Range('D30').Value = '0.09602'
$ws.Range('E30').Value = '  -4.80%  '
$ws.Range('D31').Value = '4.398'
$ws.Range('E31').Value = '  -5.16%  '
$ws.Range('D32').Value = '1.468'
$ws.Range('E32').Value = '  -4.42%  '
$ws.Range('E33').Value = '  -5.57%  '
$ws.Range('D34').Value = '0.04655'
$ws.Range('E34').Value = '  -3.90%  '
$ws.Range('D35').Value = '0.7008'
$ws.Range('E35').Value = '  -5.28%  '
$ws.Range('D36').Value = '1.096'
$ws.Range('E36').Value = '  -3.43%  '
$ws.Range('E37').Value = '  -0.70%  '
$ws.Range('E38').Value = '  -5.49%  '
$ws.Range('D39').Value = '6.286'
$ws.Range('E39').Value = '  -8.49%  '
$ws.Range('D40').Value = '2.518'
$ws.Range('E40').Value = '  -4.46%  '
$ws.Range('D41').Value = '72.61'
$ws.Range('E41').Value = '  -5.27%  '
$ws.Range('D42').Value = '0.8512'
$ws.Range('E42').Value = '  -2.46%  '
$ws.Range('D43').Value = '1.920'
$ws.Range('E43').Value = '  -5.83%  '
$ws.Range('D44').Value = '0.9994'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = '0.4155'
$ws.Range('E45').Value = '  -5.19%  '
$ws.Range('D46').Value = '103.32'
$ws.Range('E46').Value = '  -2.73%  '
$ws.Range('D47').Value = '989.78'
$ws.Range('E47').Value = '  -2.90%  '
$ws.Range('D48').Value = '7.175'
$ws.Range('E48').Value = '  -5.40%  '
$ws.Range('D49').Value = '9.193'
$ws.Range('E49').Value = '  -0.57%  '
$ws.Range('D50').Value = '34.14'
$ws.Range('E50').Value = '  -3.18%  '
$ws.Range('E51').Value = '  -5.92%  '
